$wb = $excel.ActiveWorkbook

# Rename the "CIFAR-10" dataset label to "CIFAR" in the header row of every
# table block (rows 2, 9 and 18) on both worksheets.
foreach ($ws in $wb.Worksheets) {
    foreach ($addr in @("C2", "C9", "C18")) {
        $cell = $ws.Range($addr)
        if ($cell.Value() -eq "CIFAR-10") {
            $cell.Value = "CIFAR"
        }
    }
}

# Restore the cursor/selection state recorded in the saved file: BPnumber's
# selection sits on C18 while BPlocation (the active tab) ends on C16.
$ws1 = $wb.Worksheets.Item("BPnumber")
[void]$ws1.Activate()
[void]$ws1.Range("C18").Select()

$ws2 = $wb.Worksheets.Item("BPlocation")
[void]$ws2.Activate()
[void]$ws2.Range("C16").Select()
